$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 2 (rows 7-10): "# of Users" with Total Actv. / No. of Scores / Teachers / Students ---
$ws.Range("A7").Value = "# of Users"
$ws.Range("A7:D7").Merge()

$ws.Range("A9").Value = "Total Actv."
$ws.Range("B9").Value = "No. of Scores"
$ws.Range("C9").Value = "Teachers"
$ws.Range("D9").Value = "Students"

$ws.Range("A10").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 2

# --- Spacer row 13 merged A13:F13 ---
$ws.Range("A13:F13").Merge()

# --- Table 3 (rows 15-18): "Total Statistics" ---
$ws.Range("A15").Value = "Total Statistics"

$ws.Range("A17").Value = "Total Actv."
$ws.Range("B17").Value = "No. of Scores"
$ws.Range("C17").Value = "Total Users"
$ws.Range("D17").Value = "Administrators"
$ws.Range("E17").Value = "Teachers"
$ws.Range("F17").Value = "Students"

$ws.Range("A18").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 2
